# Added live demo slide
#
# The author duplicated the final "Questions?" slide (which uses the
# "Conclusion" layout with a centered title + empty body placeholder),
# moved the duplicate so it sits right after "Mistakes Made/Lessons
# Learned" (the picture-heavy slide that was slide 10), and retitled it
# "Live Demo", leaving its body placeholder empty.

$p = $ppt.ActivePresentation

# "Questions?" is the last slide of the deck.
$questionsSlide = $p.Slides.Item($p.Slides.Count)

# Duplicate() returns a SlideRange; grab the single new Slide from it.
$liveDemo = $questionsSlide.Duplicate().Item(1)

# Re-title the duplicate (it inherited "Questions?" from the source slide).
$liveDemo.Shapes.Item(1).TextFrame.TextRange.Text = "Live Demo"

# Move it so it becomes slide 11 (right after "Mistakes Made/Lessons
# Learned", before the rest of the "Mistakes Made" / "Post-Game Press
# Conference" / "Questions?" run).
$liveDemo.MoveTo(11)
